$d = $word.ActiveDocument

# Locate the unique "Alternative Scenario A" paragraph, which follows the
# paragraph we need to insert the new separator paragraph before.
$rng = $d.Content
$found = $rng.Find.Execute("Alternative Scenario A", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertParagraphBefore()

# The freshly-inserted (empty) paragraph now sits right before "Alternative
# Scenario A". Give it the exact OOXML content/formatting of the separator
# line (including the paragraph-mark rFonts eastAsia hint).
$sepPara = $rng.Paragraphs.Item(1)
$sepXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
  '<w:t>-----------------------------------------------------------------------------------------------------</w:t></w:r>' + `
  '</w:p>'
$sepPara.Range.InsertXML($sepXml)

Write-Output "inserted separator paragraph"
